# Update vm_pu results: case re-run with 380 kV slack bus voltage (1.02 pu
# instead of 1.05 pu). Rewrites columns B-F and I-N for data rows 2-25 on
# Sheet1 with the new per-unit voltage magnitudes produced by the rerun.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters present in the data block, in the same left-to-right order
# as the values inside each row's array below (column H has no data and is
# skipped, matching the sheet's layout).
$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

# Each inner array holds the new values for row (2..25), one value per
# column in $columns, in order.
$newValues = @(
    @(1.02, 1.041958160785839, 1.051097604992053, 1.040017988306788, 1.057415273793603, 1.027628946937137, 1.047036566779596, 1.05385005141357, 1.042801528972241, 1.060150306127538, 1.019180938722228),
    @(1.02, 1.044073576991271, 1.053135948741015, 1.04185683168851, 1.059459085152227, 1.027664679818867, 1.048792430403403, 1.055697822168346, 1.04444800394525, 1.062004827721543, 1.019797561336916),
    @(1.02, 1.045437368906899, 1.054449995362226, 1.043042283862084, 1.060775930630987, 1.027684935387486, 1.049923478576926, 1.056888180187257, 1.04550856883493, 1.063198813010238, 1.020193935201068),
    @(1.02, 1.046009533595869, 1.055001274963467, 1.043539618239038, 1.061328209768736, 1.0276927655657, 1.050397771525939, 1.057387369488205, 1.04595330040775, 1.063699349618831, 1.020359951564805),
    @(1.02, 1.046105534600325, 1.055093770836404, 1.043623063186663, 1.061420863077723, 1.027694040127395, 1.050477337751748, 1.057471113746126, 1.046027907193281, 1.063783309675214, 1.020387790379049),
    @(1.02, 1.045445018770858, 1.054457366052258, 1.043048933281786, 1.060783315372962, 1.027685042706064, 1.049929820783716, 1.056894855206543, 1.045514515769146, 1.063205506730383, 1.020196155941356),
    @(1.02, 1.042674133550388, 1.051787503351642, 1.040640360279395, 1.058107172880884, 1.02764161698663, 1.047631042552186, 1.054475621319187, 1.043358973527239, 1.060778310812666, 1.01938987734905),
    @(1.02, 1.037751677299642, 1.047044127079332, 1.036361331419253, 1.053347077865635, 1.027543104352477, 1.04354003850446, 1.050171093447045, 1.039522730937773, 1.056454086610741, 1.017948640652737),
    @(1.02, 1.034441409668856, 1.04385410599897, 1.033483689148319, 1.050142119502931, 1.027462582757605, 1.040784075196985, 1.047271875499987, 1.036938302361788, 1.053537939405038, 1.016973504750287),
    @(1.02, 1.033000816257124, 1.042465816694244, 1.032231370611268, 1.048746464681396, 1.027424179872915, 1.03958357236153, 1.046009114740423, 1.035812501560061, 1.052266949372225, 1.016547742060902),
    @(1.02, 1.032464593665958, 1.04194906047765, 1.031765229927689, 1.04822683785303, 1.027409382618869, 1.039136547133105, 1.045538928858183, 1.03539329012047, 1.051793571893701, 1.016389054868521),
    @(1.02, 1.032579666599452, 1.042059956029779, 1.031865263245098, 1.048338355256114, 1.027412580809394, 1.039232486060701, 1.045639837439021, 1.03548325986992, 1.051895171183832, 1.016423118419746),
    @(1.02, 1.0329565150333, 1.042423123769575, 1.03219285931943, 1.048703537218866, 1.027422967602713, 1.039546643836766, 1.04597027247497, 1.035777870770669, 1.05222784603902, 1.01653463600634),
    @(1.02, 1.033188554005845, 1.042646738773484, 1.032394572118477, 1.048928375402331, 1.027429296616665, 1.039740059465444, 1.04617371235203, 1.035959251749436, 1.05243264807456, 1.016603273794933),
    @(1.02, 1.034536861772508, 1.043946092069931, 1.033566666512691, 1.050234575498458, 1.027465056802834, 1.040863595334561, 1.047355522601719, 1.037012873930372, 1.053622113733429, 1.017001686089629),
    @(1.02, 1.035380659530971, 1.044759247250974, 1.034300187707676, 1.051051785276018, 1.027486540433762, 1.04156642343965, 1.048094842965343, 1.037671962626534, 1.054365995612413, 1.017250649002791),
    @(1.02, 1.035872137807931, 1.045232874866954, 1.034727434307427, 1.051527691259035, 1.027498730348392, 1.041975683147666, 1.048525366092448, 1.038055750873549, 1.054799092663844, 1.017395525733584),
    @(1.02, 1.036039602464634, 1.045394256666047, 1.034873012832455, 1.051689835227964, 1.02750282897492, 1.042115114342625, 1.048672043881813, 1.038186503850062, 1.054946633379013, 1.017444867756544),
    @(1.02, 1.035290200133377, 1.044672073150759, 1.034221550507517, 1.05096418511364, 1.027484270737497, 1.041491088020983, 1.048015594545441, 1.03760131577718, 1.054286266754894, 1.017223972793035),
    @(1.02, 1.03284557383243, 1.042316210112464, 1.032096417568314, 1.048596034141595, 1.027419923667897, 1.039454162952895, 1.045872999326649, 1.035691144112139, 1.052129917041822, 1.016501811849899),
    @(1.02, 1.03130202986684, 1.040828698125802, 1.030754612567279, 1.047100016786699, 1.027376383270204, 1.038167058543649, 1.04451925118196, 1.034484118242813, 1.050766740611898, 1.016044631015928),
    @(1.02, 1.032120921175794, 1.041617863752599, 1.031466474534284, 1.047893765115179, 1.027399757555554, 1.038849994717449, 1.045237536230287, 1.035124566201603, 1.051490097684546, 1.016287291519866),
    @(1.02, 1.035331076986548, 1.044711465449675, 1.034257085138341, 1.051003770199374, 1.027485297369633, 1.041525130978807, 1.048051405687731, 1.037633240045754, 1.054322295254302, 1.017236027671483),
    @(1.02, 1.039029157652932, 1.048275170050304, 1.037471846658491, 1.054583111957935, 1.027571184376755, 1.04460259121901, 1.051289002719167, 1.040519131056116, 1.057577758147567, 1.018323716631685)
)

$startRow = 2
for ($i = 0; $i -lt $newValues.Count; $i++) {
    $rowValues = $newValues[$i]
    $row = $startRow + $i
    for ($j = 0; $j -lt $columns.Count; $j++) {
        $col = $columns[$j]
        $ws.Range("$col$row").Value = $rowValues[$j]
    }
}

Write-Output "Updated $($newValues.Count) rows x $($columns.Count) columns on $($ws.Name)."
